# Refactor code to save results in a specified folder
# Updated "IPC PO" (column C) values for rows 2-51 following the sliding
# window results re-run; DELTA (D), DELTA^2 (E), TOTAL (row 52) and MSE
# (row 53) are derived values that are recomputed from the new inputs.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newC = @(
    28.42386611449396,
    28.15315545840144,
    28.9379931344628,
    29.79311310934098,
    29.74720110808241,
    30.37403994010583,
    29.87681417992178,
    29.8686479616212,
    29.7351860059534,
    29.51825535743507,
    30.16314046967677,
    30.3018377452168,
    30.36586343673891,
    31.13742381630315,
    30.69108004341203,
    31.60766271517742,
    31.07924284054112,
    31.58446960331086,
    31.78522768947801,
    32.34002021415262,
    31.82018021041069,
    32.2387495864978,
    31.7615736532619,
    32.73824665571667,
    32.78377264559467,
    32.29237012431677,
    33.76057012806505,
    32.59403198917764,
    32.85302398170048,
    33.44428315374984,
    33.75897078536011,
    34.43429539726459,
    34.36562791655631,
    35.03381626284589,
    35.0210822963886,
    35.37734670561872,
    35.65105844094344,
    36.02596915042613,
    36.70783615399061,
    38.19650375236694,
    38.50634900725598,
    38.71439974336509,
    38.99983595173688,
    39.98974513853209,
    40.2765666450357,
    40.50062831565015,
    40.92595013331834,
    41.97935122073537,
    41.28656690707493,
    41.0054191527864
)

$startRow = 2
for ($i = 0; $i -lt $newC.Length; $i++) {
    $row = $startRow + $i
    $cVal = $newC[$i]
    $bVal = $ws.Cells.Item($row, 2).Value2
    $dVal = $cVal - $bVal
    $eVal = $dVal * $dVal

    $ws.Cells.Item($row, 3).Value2 = $cVal
    $ws.Cells.Item($row, 4).Value2 = $dVal
    $ws.Cells.Item($row, 5).Value2 = $eVal
}

$endRow = $startRow + $newC.Length - 1

$cTotal = 0
$eTotal = 0
for ($row = $startRow; $row -le $endRow; $row++) {
    $cTotal = $cTotal + $ws.Cells.Item($row, 4).Value2
    $eTotal = $eTotal + $ws.Cells.Item($row, 5).Value2
}

$ws.Cells.Item(52, 3).Value2 = $cTotal
$ws.Cells.Item(52, 5).Value2 = $eTotal
$ws.Cells.Item(53, 5).Value2 = $eTotal / $newC.Length
